# Strategy AI Agent 45-Day Plan - mark Week 1 tasks (Days 2-7, rows 3-8) as
# "Done" in the Status column, matching the already-completed Day 1 (row 2),
# and move the active selection to the next open task (E9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3:E8").Value = "Done"

$ws.Range("E9").Select()
